# Update Sheets via scheduled runner: refresh market price / profit figures
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 391.85715
$ws.Range("I2").Value = 430
$ws.Range("J2").Value = 363.25
$ws.Range("K2").Value = 430
$ws.Range("L2").Value = 363.25
$ws.Range("M2").Value = -317
$ws.Range("N2").Value = -589.25
$ws.Range("H28").Value = 1294.7
$ws.Range("I28").Value = 1243.1154
$ws.Range("K28").Value = 1243.1154
$ws.Range("M28").Value = -758.1153999999999
$ws.Range("H100").Value = 5848.8096
$ws.Range("I100").Value = 3502.5
$ws.Range("K100").Value = 3502.5
$ws.Range("M100").Value = -2961.5
$ws.Range("H132").Value = 735.0364
$ws.Range("I132").Value = 730.6981
$ws.Range("K132").Value = 2192.0943
$ws.Range("M132").Value = 337.9057000000003
$ws.Range("H138").Value = 4359.7407
$ws.Range("J138").Value = 4359.7407
$ws.Range("L138").Value = 13079.2221
$ws.Range("N138").Value = -23359.2221

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15317.515
$ws.Range("I32").Value = 11620
$ws.Range("J32").Value = 39556.777
$ws.Range("K32").Value = 11620
$ws.Range("L32").Value = 39556.777
$ws.Range("M32").Value = -11333
$ws.Range("N32").Value = -40130.777
$ws.Range("H61").Value = 4513.2
$ws.Range("I61").Value = 2585.5557
$ws.Range("K61").Value = 2585.5557
$ws.Range("M61").Value = -2373.5557
$ws.Range("H102").Value = 55557016
$ws.Range("I102").Value = 1749.8
$ws.Range("K102").Value = 1749.8
$ws.Range("M102").Value = -127.8
$ws.Range("H110").Value = 2598.6453
$ws.Range("I110").Value = 2537.7856
$ws.Range("K110").Value = 2537.7856
$ws.Range("M110").Value = -492.7856000000002
$ws.Range("H136").Value = 4513.2
$ws.Range("I136").Value = 2585.5557
$ws.Range("K136").Value = 7756.6671
$ws.Range("M136").Value = -5206.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 49999.168
$ws.Range("J70").Value = 49999.168
$ws.Range("L70").Value = 49999.168
$ws.Range("N70").Value = -50585.168
$ws.Range("H73").Value = 49999.168
$ws.Range("J73").Value = 49999.168
$ws.Range("L73").Value = 49999.168
$ws.Range("N73").Value = -52027.168
$ws.Range("H135").Value = 66555.55499999999
$ws.Range("I135").Value = 39000
$ws.Range("K135").Value = 39000
$ws.Range("M135").Value = -33930

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 593
$ws.Range("I16").Value = 492
$ws.Range("K16").Value = 492
$ws.Range("M16").Value = -205
$ws.Range("H22").Value = 754.3
$ws.Range("I22").Value = 257.66666
$ws.Range("K22").Value = 257.66666
$ws.Range("M22").Value = 92.33334000000002
$ws.Range("H33").Value = 2500.5
$ws.Range("I33").Value = 2500.5
$ws.Range("K33").Value = 2500.5
$ws.Range("M33").Value = -2121.5
$ws.Range("H58").Value = 2366.8333
$ws.Range("I58").Value = 1200.2858
$ws.Range("K58").Value = 1200.2858
$ws.Range("M58").Value = -997.2858000000001
$ws.Range("H94").Value = 6419.65
$ws.Range("I94").Value = 15448.429
$ws.Range("K94").Value = 15448.429
$ws.Range("M94").Value = -14997.429
$ws.Range("H109").Value = 12666.25
$ws.Range("J109").Value = 12666.25
$ws.Range("L109").Value = 12666.25
$ws.Range("N109").Value = -14746.25
$ws.Range("H113").Value = 593
$ws.Range("I113").Value = 492
$ws.Range("K113").Value = 492
$ws.Range("M113").Value = 1678
$ws.Range("H136").Value = 2366.8333
$ws.Range("I136").Value = 1200.2858
$ws.Range("K136").Value = 3600.8574
$ws.Range("M136").Value = -1050.8574
$ws.Range("H141").Value = 362634.84
$ws.Range("J141").Value = 362634.84
$ws.Range("L141").Value = 362634.84
$ws.Range("N141").Value = -372994.84

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3878
$ws.Range("I113").Value = 2769.3635
$ws.Range("J113").Value = 7943
$ws.Range("K113").Value = 2769.3635
$ws.Range("L113").Value = 7943
$ws.Range("M113").Value = -599.3634999999999
$ws.Range("N113").Value = -12283
$ws.Range("H126").Value = 3046.3333
$ws.Range("J126").Value = 7998.3335
$ws.Range("L126").Value = 23995.0005
$ws.Range("N126").Value = -28935.0005
$ws.Range("H132").Value = 4586.8076
$ws.Range("I132").Value = 4170.28
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 12510.84
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -9980.84
$ws.Range("N132").Value = -50060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3916.4614
$ws.Range("I7").Value = 2431.6
$ws.Range("J7").Value = 5941.273
$ws.Range("K7").Value = 2431.6
$ws.Range("L7").Value = 5941.273
$ws.Range("M7").Value = -2319.6
$ws.Range("N7").Value = -6165.273
$ws.Range("H22").Value = 1329.4445
$ws.Range("I22").Value = 1177.5834
$ws.Range("K22").Value = 1177.5834
$ws.Range("M22").Value = -882.5834
$ws.Range("H27").Value = 1329.4445
$ws.Range("I27").Value = 1177.5834
$ws.Range("K27").Value = 1177.5834
$ws.Range("M27").Value = -1070.5834
$ws.Range("H61").Value = 2680.0312
$ws.Range("I61").Value = 1715.826
$ws.Range("K61").Value = 1715.826
$ws.Range("M61").Value = -1513.826
$ws.Range("H80").Value = 91851.57000000001
$ws.Range("J80").Value = 91851.57000000001
$ws.Range("L80").Value = 91851.57000000001
$ws.Range("N80").Value = -94097.57000000001
$ws.Range("H83").Value = 91851.57000000001
$ws.Range("J83").Value = 91851.57000000001
$ws.Range("L83").Value = 275554.71
$ws.Range("N83").Value = -286786.71
$ws.Range("H93").Value = 21670316
$ws.Range("I93").Value = 4379.3
$ws.Range("J93").Value = 130000000
$ws.Range("K93").Value = 4379.3
$ws.Range("L93").Value = 130000000
$ws.Range("M93").Value = -3131.3
$ws.Range("N93").Value = -130002496
$ws.Range("H113").Value = 2680.0312
$ws.Range("I113").Value = 1715.826
$ws.Range("K113").Value = 1715.826
$ws.Range("M113").Value = 454.174
$ws.Range("H126").Value = 3916.4614
$ws.Range("I126").Value = 2431.6
$ws.Range("J126").Value = 5941.273
$ws.Range("K126").Value = 7294.799999999999
$ws.Range("L126").Value = 17823.819
$ws.Range("M126").Value = -4824.799999999999
$ws.Range("N126").Value = -22763.819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 34000
$ws.Range("J99").Value = 34000
$ws.Range("L99").Value = 34000
$ws.Range("M99").Value = -39990
$ws.Range("H107").Value = 438.6875
$ws.Range("I107").Value = 372.7857
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 1118.3571
$ws.Range("L107").Value = 2700
$ws.Range("M107").Value = 801.6428999999998
$ws.Range("N107").Value = -6540
$ws.Range("H109").Value = 28667.285
$ws.Range("J109").Value = 28461.54
$ws.Range("L109").Value = 28461.54
$ws.Range("N109").Value = -31235.54
$ws.Range("H113").Value = 511.8
$ws.Range("I113").Value = 520.9545000000001
$ws.Range("K113").Value = 1562.8635
$ws.Range("M113").Value = 607.1364999999998
$ws.Range("H126").Value = 1451.875
$ws.Range("I126").Value = 1215.3334
$ws.Range("K126").Value = 3646.0002
$ws.Range("M126").Value = -1176.0002

Write-Host "Updated 180 cells across 7 sheets"
